# Corrected typo in sprint 2 burndown chart:
# the worksheet was mistakenly named "Sprint 3" instead of "Sprint 2".
# Renaming the sheet does not automatically retarget the embedded
# chart's series formulas in this host, so update those explicitly
# as well (mirrors what Excel itself would do when a referenced sheet
# is renamed).

$wb = $excel.ActiveWorkbook
$oldName = "Sprint 3"
$newName = "Sprint 2"

$ws = $wb.Worksheets.Item($oldName)
$ws.Name = $newName

foreach ($co in $ws.ChartObjects()) {
    $chart = $co.Chart
    foreach ($ser in $chart.SeriesCollection()) {
        $ser.Formula = $ser.Formula.Replace("'" + $oldName + "'", "'" + $newName + "'")
    }
}
